$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "purchasing location" notes next to a few parts, per the commit message
# ("Listed locations for purchasing parts").
$ws.Range("I11").Value = "Out of stock: Substituting 587-2985-1-ND"
$ws.Range("I15").Value = "signpost"
$ws.Range("I17").Value = "stock"

# Matches the author's final cursor position recorded in the sheet view.
[void]$ws.Range("K23").Select()

# The sheet was re-flowed (columns nudged to their new optimal widths) after
# the notes above were typed in - column E especially grows to fit the
# longest "Description" text now that the sheet is wider.
$ws.Columns.Item(1).ColumnWidth = 3.33
$ws.Columns.Item(2).Resize(1, 2).ColumnWidth = 18.67
$ws.Columns.Item(4).ColumnWidth = 17.83
$ws.Columns.Item(5).ColumnWidth = 118.17
$ws.Columns.Item(6).ColumnWidth = 40.83
$ws.Columns.Item(7).ColumnWidth = 25.17
